$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for January 2025 above the current first data row (row 6) ---
$ws.Rows.Item(6).Insert()

# --- Populate the new row's values ---
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Ene."
$ws.Range("D6").Value = 3715.108

# --- Resize the table / AutoFilter range to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B5:D90"))

# --- Update the "last updated" caption ---
$ws.Range("B91").Value = "Actualización: Enero 2025."

# --- Formatting for the new row: center aligned, thin borders (top+bottom, plus left on B) ---
$b6 = $ws.Range("B6")
$c6 = $ws.Range("C6")
$d6 = $ws.Range("D6")

$b6.HorizontalAlignment = -4108
$c6.HorizontalAlignment = -4108

$d6.NumberFormat = "#,##0.0"

$rowRng = $ws.Range("B6:D6")
$rowRng.Borders.Item(8).LineStyle = 1
$rowRng.Borders.Item(9).LineStyle = 1
$b6.Borders.Item(7).LineStyle = 1

Write-Host "done"
